{"js": "// Update maze clue to reference coordinate system.\n// Three edits inside word/document.xml:\n//  1) \"A map ... front door.\" paragraph -> reworded + centered.\n//  2) \"Start from the top left ... fill in 0's with black and 1's with white\"\n//     -> \"... and create a BitPic of the original map.\"\n//  3) \"Hint: Type pwd ...\" paragraph -> completely reworded hint about the\n//     coordinate system (with \"coordinate system\" bolded).\n\nconst body = context.document.body;\n\n// --- Edit 1: reword the intro paragraph and center it ---------------------\nconst intro = body.search(\n  \"A map of how to get through the maze (starting on the left) and get to the Hogwarts front door. Unfortunately, the \",\n  { matchCase: true }\n);\nintro.load(\"items\");\nawait context.sync();\n\nif (intro.items.length > 0) {\n  const introRange = intro.items[0];\n  const introPara = introRange.paragraphs.getFirst();\n  introPara.alignment = Word.Alignment.centered;\n  introRange.insertText(\n    \"Below is a map of how to get through the maze (starting on the left) and get to the Hogwarts.  Unfortunately, the \",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// --- Edit 2: swap the bitmap instructions for the BitPic instructions -----\nconst bitmapTail = body.search(\n  \" and fill in 0\\u2019s with black and 1\\u2019s with white\",\n  { matchCase: true }\n);\nbitmapTail.load(\"items\");\nawait context.sync();\n\nif (bitmapTail.items.length > 0) {\n  bitmapTail.items[0].insertText(\n    \" and create a BitPic of the original map.\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// --- Edit 3: replace the whole hint paragraph ------------------------------\n// The original paragraph contains several w:proofErr spell-check markers\n// (around \"pwd\", \"apparate\", \"maze_entrance\"). A plain insertText(replace)\n// would leave those orphaned, so clear() the paragraph first to drop all of\n// its run/markup children, then type the new sentence back in.\nconst hintSearch = body.search(\n  \"Hint: Type pwd before going through the maze so you know the full \\u201crecipe\\u201d of the entrance. This will be helpful if you get lost and need to apparate back to the maze_entrance.\",\n  { matchCase: true }\n);\nhintSearch.load(\"items\");\nawait context.sync();\n\nif (hintSearch.items.length > 0) {\n  const hintPara = hintSearch.items[0].paragraphs.getFirst();\n  hintPara.clear();\n  await context.sync();\n\n  hintPara.font.name = \"Garamond\";\n  hintPara.insertText(\n    \"Hint: Notice that the map has a coordinate system. Be sure to use it as you try to find your way through the maze.\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// Bold just \"coordinate system\" within the new hint text.\nconst coordPhrase = body.search(\"coordinate system\", { matchCase: true });\ncoordPhrase.load(\"items\");\nawait context.sync();\nif (coordPhrase.items.length > 0) {\n  coordPhrase.items[0].font.bold = true;\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Edit 1: reword the intro paragraph and center it ----------------------\n$r1 = $d.Content\n$r1.Find.ClearFormatting()\n$r1.Find.Text = \"A map of how to get through the maze (starting on the left) and get to the Hogwarts front door. Unfortunately, the \"\n$r1.Find.Replacement.ClearFormatting()\n$r1.Find.Replacement.Text = \"Below is a map of how to get through the maze (starting on the left) and get to the Hogwarts.  Unfortunately, the \"\n$r1.Find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$r1c = $d.Content\n$r1c.Find.ClearFormatting()\n$r1c.Find.Text = \"Below is a map\"\n$r1c.Find.Execute() | Out-Null\n$r1c.Paragraphs(1).Alignment = 1\n\n# --- Edit 2: swap the bitmap instructions for the BitPic instructions ------\n$r2 = $d.Content\n$r2.Find.ClearFormatting()\n$r2.Find.Text = \" and fill in 0\u2019s with black and 1\u2019s with white\"\n$r2.Find.Replacement.ClearFormatting()\n$r2.Find.Replacement.Text = \" and create a BitPic of the original map.\"\n$r2.Find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n# --- Edit 3: replace the whole hint paragraph -------------------------------\n$r3 = $d.Content\n$r3.Find.ClearFormatting()\n$r3.Find.Text = \"Hint: Type pwd before going through the maze so you know the full \u201crecipe\u201d of the entrance. This will be helpful if you get lost and need to apparate back to the maze_entrance.\"\n$r3.Find.Replacement.ClearFormatting()\n$r3.Find.Replacement.Text = \"Hint: Notice that the map has a coordinate system. Be sure to use it as you try to find your way through the maze.\"\n$r3.Find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null\n\n$r3b = $d.Content\n$r3b.Find.ClearFormatting()\n$r3b.Find.Text = \"coordinate system\"\n$r3b.Find.Execute() | Out-Null\n$r3b.Font.Bold = 1\n"}
